# Apply cryptocurrency price/volume updates to match the Sun Aug 13 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.379.80"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.849.77"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.58"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6291"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07614"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2916"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.025"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6817"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001052"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.17"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.131"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "29.390.69"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.58"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.0000"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "158.74"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1392"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.448"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.70"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +9.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.476"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05630"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.113"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.076"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.834"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6974"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.583"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01816"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").Value = "1.235.93"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.417"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9028"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.0000"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.53"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.66"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4002"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1157"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.686"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.991"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05703"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  -0.10%  "
